# Commit: "added a feature that support persian title on import excel the"
# -> Replace the English column header strings with their Persian
#    translations, and move the sheet's active selection from C5 to C2
#    (the cell just below the headers), matching the author's workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "ستون اول"
$ws.Range("B1").Value = "ستون دوم"
$ws.Range("C1").Value = "ستون سوم"

$ws.Range("C2").Select()
